$wb = $excel.ActiveWorkbook

# ===== Sheet 1: "展览" =====
$ws1 = $wb.Worksheets.Item(1)

# Simple F-column (want-to-go count) refreshes for rows untouched by the
# new-row insertion below.
$ws1.Range("F5").Value = 69
$ws1.Range("F6").Value = 2357
$ws1.Range("F7").Value = 67
$ws1.Range("F10").Value = 70
$ws1.Range("F12").Value = 1495
$ws1.Range("F13").Value = 18
$ws1.Range("F14").Value = 590
$ws1.Range("F15").Value = 548
$ws1.Range("F16").Value = 1003
$ws1.Range("F17").Value = 479
$ws1.Range("F18").Value = 3291
$ws1.Range("F20").Value = 132
$ws1.Range("F21").Value = 3243
$ws1.Range("F22").Value = 708
$ws1.Range("F23").Value = 592

# Two brand-new listings were added around 2024-11-16/17, pushing the
# existing tail of the table down. Insert two blank rows at the right spots.
$ws1.Rows.Item(24).Insert()
$ws1.Rows.Item(26).Insert()

# Re-apply the index-column formatting (bold, centered, thin border) that
# the other rows in column A carry, since a freshly inserted row starts blank.
foreach ($r in @(24, 26)) {
    $idxCell = $ws1.Range("A" + $r)
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1
}

# Write final content for rows 24-31 (both new rows and the shifted-down
# existing rows, some of which also got a refreshed F-column count).
# -- row 24 --
$ws1.Range("A24").Value = 23
$ws1.Range("B24").NumberFormat = "@"
$ws1.Range("B24").Value = '2024-11-16'
$ws1.Range("C24").NumberFormat = "@"
$ws1.Range("C24").Value = '北京·万游引力国潮动漫嘉年华S9 知名配音演员 桑毓泽 现场签售票'
$ws1.Range("D24").NumberFormat = "@"
$ws1.Range("D24").Value = '金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心'
$ws1.Range("E24").NumberFormat = "@"
$ws1.Range("E24").Value = '2024.11.16 11:00-11.16 17:00'
$ws1.Range("F24").Value = 0
$ws1.Range("G24").Value = 52
$ws1.Range("H24").NumberFormat = "@"
$ws1.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=93151'
$ws1.Range("I24").NumberFormat = "@"
$ws1.Range("I24").Value = '//i2.hdslb.com/bfs/openplatform/202410/lACYIuZh1728312527225.jpeg'

# -- row 25 --
$ws1.Range("A25").Value = 24
$ws1.Range("B25").NumberFormat = "@"
$ws1.Range("B25").Value = '2024-11-16'
$ws1.Range("C25").NumberFormat = "@"
$ws1.Range("C25").Value = '北京·原神×星穹铁道×绝区零同人ONLY'
$ws1.Range("D25").NumberFormat = "@"
$ws1.Range("D25").Value = '永外高庄138号 北京大红门国际会展中心'
$ws1.Range("E25").NumberFormat = "@"
$ws1.Range("E25").Value = '2024.11.16 10:00-11.16 17:00'
$ws1.Range("F25").Value = 268
$ws1.Range("G25").Value = 60
$ws1.Range("H25").NumberFormat = "@"
$ws1.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=92358'
$ws1.Range("I25").NumberFormat = "@"
$ws1.Range("I25").Value = '//i1.hdslb.com/bfs/openplatform/202409/A24uEchR1726118358020.jpeg'

# -- row 26 --
$ws1.Range("A26").Value = 25
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("B26").Value = '2024-11-17'
$ws1.Range("C26").NumberFormat = "@"
$ws1.Range("C26").Value = '北京·万游引力国潮动漫嘉年华S9 知名唱见 茶师 现场签售券'
$ws1.Range("D26").NumberFormat = "@"
$ws1.Range("D26").Value = '金蝉西路甲1号 北京酷车国际汇展中心'
$ws1.Range("E26").NumberFormat = "@"
$ws1.Range("E26").Value = '2024.11.17 11:00-11.17 17:00'
$ws1.Range("F26").Value = 0
$ws1.Range("G26").Value = 52
$ws1.Range("H26").NumberFormat = "@"
$ws1.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=93150'
$ws1.Range("I26").NumberFormat = "@"
$ws1.Range("I26").Value = '//i0.hdslb.com/bfs/openplatform/202410/fpaNLvw11728312099093.jpeg'

# -- row 27 --
$ws1.Range("A27").Value = 26
$ws1.Range("B27").NumberFormat = "@"
$ws1.Range("B27").Value = '2024-11-23'
$ws1.Range("C27").NumberFormat = "@"
$ws1.Range("C27").Value = '北京·代号鸢only同人展（如鸢同人only）'
$ws1.Range("D27").NumberFormat = "@"
$ws1.Range("D27").Value = '北花园路1号超级蜂巢C座 超级蜂巢国际会议中心'
$ws1.Range("E27").NumberFormat = "@"
$ws1.Range("E27").Value = '2024.11.23 10:00-11.23 17:00'
$ws1.Range("F27").Value = 1073
$ws1.Range("G27").Value = 68
$ws1.Range("H27").NumberFormat = "@"
$ws1.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=90673'
$ws1.Range("I27").NumberFormat = "@"
$ws1.Range("I27").Value = '//i1.hdslb.com/bfs/openplatform/202409/p7GKpOhb1727406520811.jpeg'

# -- row 28 --
$ws1.Range("A28").Value = 27
$ws1.Range("B28").NumberFormat = "@"
$ws1.Range("B28").Value = '2024-11-30'
$ws1.Range("C28").NumberFormat = "@"
$ws1.Range("C28").Value = '北京·蔚蓝档案only同人展'
$ws1.Range("D28").NumberFormat = "@"
$ws1.Range("D28").Value = '北花园路1号超级蜂巢C座 超级蜂巢国际会议中心'
$ws1.Range("E28").NumberFormat = "@"
$ws1.Range("E28").Value = '2024.11.30 10:00-11.30 17:00'
$ws1.Range("F28").Value = 752
$ws1.Range("G28").Value = 68
$ws1.Range("H28").NumberFormat = "@"
$ws1.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=92109'
$ws1.Range("I28").NumberFormat = "@"
$ws1.Range("I28").Value = '//i1.hdslb.com/bfs/openplatform/202409/rG5Ps2Em1727063078808.jpeg'

# -- row 29 --
$ws1.Range("A29").Value = 28
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = '2024-12-14'
$ws1.Range("C29").NumberFormat = "@"
$ws1.Range("C29").Value = '北京·奇想派对第五届'
$ws1.Range("D29").NumberFormat = "@"
$ws1.Range("D29").Value = '学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里'
$ws1.Range("E29").NumberFormat = "@"
$ws1.Range("E29").Value = '2024.12.14 10:00-12.15 17:30'
$ws1.Range("F29").Value = 43
$ws1.Range("G29").Value = 45
$ws1.Range("H29").NumberFormat = "@"
$ws1.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=91077'
$ws1.Range("I29").NumberFormat = "@"
$ws1.Range("I29").Value = '//i1.hdslb.com/bfs/openplatform/202408/zMayUoC81724229782742.jpeg'

# -- row 30 --
$ws1.Range("A30").Value = 29
$ws1.Range("B30").NumberFormat = "@"
$ws1.Range("B30").Value = '2024-12-28'
$ws1.Range("C30").NumberFormat = "@"
$ws1.Range("C30").Value = '北京·第20届IJOY漫展xCGF游戏节'
$ws1.Range("D30").NumberFormat = "@"
$ws1.Range("D30").Value = '天辰东路7号 北京国家会议中心'
$ws1.Range("E30").NumberFormat = "@"
$ws1.Range("E30").Value = '2024.12.28 09:00-12.29 17:00'
$ws1.Range("F30").Value = 894
$ws1.Range("G30").Value = 8.8
$ws1.Range("H30").NumberFormat = "@"
$ws1.Range("H30").Value = 'https://show.bilibili.com/platform/detail.html?id=92633'
$ws1.Range("I30").NumberFormat = "@"
$ws1.Range("I30").Value = '//i0.hdslb.com/bfs/openplatform/202409/EQg8HwjJ1726734597607.jpeg'

# -- row 31 --
$ws1.Range("A31").Value = 30
$ws1.Range("B31").NumberFormat = "@"
$ws1.Range("B31").Value = '2025-01-17'
$ws1.Range("C31").NumberFormat = "@"
$ws1.Range("C31").Value = ' 北京·第21届IJOY漫展xCGF游戏节'
$ws1.Range("D31").NumberFormat = "@"
$ws1.Range("D31").Value = '天辰东路7号 北京国家会议中心'
$ws1.Range("E31").NumberFormat = "@"
$ws1.Range("E31").Value = '2025.01.17 09:00-01.19 17:00'
$ws1.Range("F31").Value = 863
$ws1.Range("G31").Value = 8.8
$ws1.Range("H31").NumberFormat = "@"
$ws1.Range("H31").Value = 'https://show.bilibili.com/platform/detail.html?id=92634'
$ws1.Range("I31").NumberFormat = "@"
$ws1.Range("I31").Value = '//i0.hdslb.com/bfs/openplatform/202409/ASXIizNW1726735204415.jpeg'

# ===== Sheet 2: "演出" =====
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 78
$ws2.Range("F14").Value = 193
$ws2.Range("F19").Value = 220
$ws2.Range("F20").Value = 159

# ===== Sheet 3: "本地生活" =====
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 2945
$ws3.Range("F4").Value = 382
$ws3.Range("F6").Value = 456

# ===== Sheet 4: "全部类型" =====
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 382
$ws4.Range("F11").Value = 69
$ws4.Range("F12").Value = 456
$ws4.Range("F13").Value = 2357
$ws4.Range("F14").Value = 67
$ws4.Range("F19").Value = 78
$ws4.Range("F24").Value = 1495
$ws4.Range("F25").Value = 1495
$ws4.Range("F27").Value = 18
$ws4.Range("F28").Value = 548
$ws4.Range("F30").Value = 193
$ws4.Range("F31").Value = 1004
$ws4.Range("F32").Value = 479
$ws4.Range("F34").Value = 3291
$ws4.Range("F35").Value = 132
$ws4.Range("F36").Value = 3243
$ws4.Range("F37").Value = 708
$ws4.Range("F39").Value = 592
$ws4.Range("F40").Value = 268
$ws4.Range("F41").Value = 1073
$ws4.Range("F44").Value = 220
$ws4.Range("F45").Value = 159
$ws4.Range("F49").Value = 894
$ws4.Range("F50").Value = 863

Write-Output "edit complete"
